$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D and E contain text-formatted numbers/percentages (e.g. "1.00", "3.340.76",
# "  -4.59%  "). Excel auto-converts numeric-looking input to real numbers on Value
# assignment, which would corrupt both the formatting (e.g. "1.00" -> 1) and precision
# (e.g. "570.78" -> 570.77999999999997). Force the range to Text format first so the
# literal strings are preserved exactly, then restore the original (default) style so the
# workbook formatting is left unchanged.
$priceVolumeRange = $ws.Range("D2:E51")
$priceVolumeRange.NumberFormat = "@"

$ws.Range("D2").Value = '60.455.70'
$ws.Range("E2").Value = '  -4.39%  '
$ws.Range("D3").Value = '3.342.87'
$ws.Range("E3").Value = '  -2.23%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").Value = '570.78'
$ws.Range("E5").Value = '  -1.34%  '
$ws.Range("D6").Value = '131.91'
$ws.Range("E6").Value = '  +3.38%  '
$ws.Range("D8").Value = '3.342.60'
$ws.Range("E8").Value = '  -2.22%  '
$ws.Range("D9").Value = '0.475'
$ws.Range("E9").Value = '  -0.57%  '
$ws.Range("D10").Value = '7.62'
$ws.Range("E10").Value = '  +1.68%  '
$ws.Range("E11").Value = '  -1.94%  '
$ws.Range("D12").Value = '0.383'
$ws.Range("E12").Value = '  +1.16%  '
$ws.Range("D13").Value = '3.914.28'
$ws.Range("E13").Value = '  -2.25%  '
$ws.Range("E14").Value = '  +0.41%  '
$ws.Range("D15").Value = '0.0000172'
$ws.Range("E15").Value = '  -1.81%  '
$ws.Range("D16").Value = '3.338.98'
$ws.Range("E16").Value = '  -2.44%  '
$ws.Range("E17").Value = '  -0.58%  '
$ws.Range("D18").Value = '60.562.08'
$ws.Range("D19").Value = '13.82'
$ws.Range("E19").Value = '  +4.72%  '
$ws.Range("D20").Value = '5.75'
$ws.Range("E20").Value = '  +1.02%  '
$ws.Range("D21").Value = '9.25'
$ws.Range("E21").Value = '  -4.72%  '
$ws.Range("D22").Value = '372.72'
$ws.Range("E22").Value = '  -1.90%  '
$ws.Range("E23").Value = '  +0.11%  '
$ws.Range("B24").Value = 'WrappedeETH'
$ws.Range("C24").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D24").Value = '3.476.91'
$ws.Range("E24").Value = '  -2.27%  '
$ws.Range("B25").Value = 'Dai'
$ws.Range("C25").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D25").Value = '1.00'
$ws.Range("E25").Value = '  -0.07%  '
$ws.Range("D26").Value = '69.84'
$ws.Range("E26").Value = '  -4.18%  '
$ws.Range("E27").Value = '  +4.57%  '
$ws.Range("E28").Value = '  +18.83%  '
$ws.Range("D29").Value = '7.58'
$ws.Range("E29").Value = '  +7.69%  '
$ws.Range("D30").Value = '1.00'
$ws.Range("E30").Value = '  +0.50%  '
$ws.Range("D31").Value = '8.05'
$ws.Range("E31").Value = '  +1.88%  '
$ws.Range("E32").Value = '  -1.03%  '
$ws.Range("E33").Value = '  +0.05%  '
$ws.Range("E34").Value = '  -0.04%  '
$ws.Range("D35").Value = '3.372.52'
$ws.Range("E35").Value = '  -2.23%  '
$ws.Range("D36").Value = '23.06'
$ws.Range("E36").Value = '  +0.89%  '
$ws.Range("D37").Value = '5.43'
$ws.Range("E37").Value = '  +2.01%  '
$ws.Range("D38").Value = '6.92'
$ws.Range("E38").Value = '  +2.37%  '
$ws.Range("E39").Value = '  +1.46%  '
$ws.Range("D40").Value = '161.85'
$ws.Range("E40").Value = '  -1.02%  '
$ws.Range("E41").Value = '  +2.00%  '
$ws.Range("D42").Value = '0.999'
$ws.Range("E42").Value = '  -0.07%  '
$ws.Range("D43").Value = '1.22'
$ws.Range("E43").Value = '  +10.74%  '
$ws.Range("E44").Value = '  +2.13%  '
$ws.Range("D45").Value = '41.26'
$ws.Range("E45").Value = '  -2.07%  '
$ws.Range("D46").Value = '0.749'
$ws.Range("E46").Value = '  -4.07%  '
$ws.Range("D47").Value = '23.51'
$ws.Range("E47").Value = '  +2.58%  '
$ws.Range("E48").Value = '  -0.68%  '
$ws.Range("E49").Value = '  +3.09%  '
$ws.Range("D50").Value = '22.59'
$ws.Range("E50").Value = '  +10.71%  '
$ws.Range("D51").Value = '0.895'
$ws.Range("E51").Value = '  +3.11%  '

$priceVolumeRange.Style = "Normal"
